# Changing Astro location to Dabaca
# - Shift all timestamps in column A (rows 2-97) forward by 8 days.
# - Update the "Actual Production (MW)" values in column B for rows 22-44
#   to reflect the new location's production curve.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2 through 97) by +8 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 8
}

# New production values (column B) for rows 22-44.
$newB = @{
    22 = 4
    23 = 14
    24 = 24
    25 = 47
    26 = 79
    27 = 110
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
